$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '22.412.29'
Set-TextValue $ws.Range("E2") '  -0.23%  '
Set-TextValue $ws.Range("D3") '1.566.42'
Set-TextValue $ws.Range("E3") '  -0.54%  '
Set-TextValue $ws.Range("E4") '  -0.01%  '
Set-TextValue $ws.Range("E5") '  -0.08%  '
Set-TextValue $ws.Range("D6") '284.53'
Set-TextValue $ws.Range("E6") '  -2.49%  '
Set-TextValue $ws.Range("D7") '0.3644'
Set-TextValue $ws.Range("E7") '  -2.41%  '
Set-TextValue $ws.Range("D8") '48.57'
Set-TextValue $ws.Range("E8") '  -2.74%  '
Set-TextValue $ws.Range("D9") '0.3332'
Set-TextValue $ws.Range("E9") '  -2.00%  '
Set-TextValue $ws.Range("D10") '1.124'
Set-TextValue $ws.Range("E10") '  -1.86%  '
Set-TextValue $ws.Range("D11") '0.07400'
Set-TextValue $ws.Range("E11") '  -2.37%  '
Set-TextValue $ws.Range("E12") '  +0.01%  '
Set-TextValue $ws.Range("D13") '20.73'
Set-TextValue $ws.Range("E13") '  -2.87%  '
Set-TextValue $ws.Range("D14") '5.950'
Set-TextValue $ws.Range("E14") '  -1.18%  '
Set-TextValue $ws.Range("D15") '6.900'
Set-TextValue $ws.Range("E15") '  -0.98%  '
Set-TextValue $ws.Range("D16") '1.566.48'
Set-TextValue $ws.Range("E16") '  -0.44%  '
Set-TextValue $ws.Range("D17") '0.00001104'
Set-TextValue $ws.Range("E17") '  -1.88%  '
Set-TextValue $ws.Range("D18") '88.15'
Set-TextValue $ws.Range("E18") '  -3.11%  '
Set-TextValue $ws.Range("D19") '0.06701'
Set-TextValue $ws.Range("E19") '  -0.96%  '
Set-TextValue $ws.Range("E20") '  -0.07%  '
Set-TextValue $ws.Range("D21") '6.348'
Set-TextValue $ws.Range("E21") '  +0.60%  '
Set-TextValue $ws.Range("D22") '16.17'
Set-TextValue $ws.Range("E22") '  -0.98%  '
Set-TextValue $ws.Range("D23") '12.02'
Set-TextValue $ws.Range("E23") '  -1.35%  '
Set-TextValue $ws.Range("D24") '22.413.66'
Set-TextValue $ws.Range("E24") '  -0.20%  '
Set-TextValue $ws.Range("D25") '2.383'
Set-TextValue $ws.Range("E25") '  +1.58%  '
Set-TextValue $ws.Range("D26") '2.537'
Set-TextValue $ws.Range("E26") '  -5.63%  '
Set-TextValue $ws.Range("D27") '150.15'
Set-TextValue $ws.Range("E27") '  +1.02%  '
Set-TextValue $ws.Range("D28") '19.35'
Set-TextValue $ws.Range("E28") '  -3.90%  '
Set-TextValue $ws.Range("D29") '5.008'
Set-TextValue $ws.Range("E29") '  -0.05%  '
Set-TextValue $ws.Range("D30") '123.75'
Set-TextValue $ws.Range("E30") '  -1.48%  '
Set-TextValue $ws.Range("D31") '1.744.37'
Set-TextValue $ws.Range("E31") '  -0.20%  '
Set-TextValue $ws.Range("D32") '1.055'
Set-TextValue $ws.Range("E32") '  -0.04%  '
Set-TextValue $ws.Range("D33") '6.101'
Set-TextValue $ws.Range("E33") '  -1.86%  '
Set-TextValue $ws.Range("D34") '1.991'
Set-TextValue $ws.Range("E34") '  +0.29%  '
Set-TextValue $ws.Range("D35") '9.800'
Set-TextValue $ws.Range("E35") '  -0.40%  '
Set-TextValue $ws.Range("D36") '0.08268'
Set-TextValue $ws.Range("E36") '  -1.57%  '
Set-TextValue $ws.Range("D37") '0.02413'
Set-TextValue $ws.Range("E37") '  -3.07%  '
Set-TextValue $ws.Range("D38") '0.2232'
Set-TextValue $ws.Range("E38") '  -2.75%  '
Set-TextValue $ws.Range("D39") '0.06406'
Set-TextValue $ws.Range("E39") '  -1.97%  '
Set-TextValue $ws.Range("D40") '5.371'
Set-TextValue $ws.Range("E40") '  -2.07%  '
Set-TextValue $ws.Range("D41") '1.282'
Set-TextValue $ws.Range("E41") '  -6.77%  '
Set-TextValue $ws.Range("E42") '  +0.29%  '
Set-TextValue $ws.Range("D43") '11.13'
Set-TextValue $ws.Range("E43") '  -1.77%  '
Set-TextValue $ws.Range("D44") '1.001'
Set-TextValue $ws.Range("E44") '  -0.08%  '
Set-TextValue $ws.Range("D45") '13.75'
Set-TextValue $ws.Range("E45") '  -1.83%  '
Set-TextValue $ws.Range("D46") '0.6053'
Set-TextValue $ws.Range("E46") '  +3.88%  '
Set-TextValue $ws.Range("D47") '3.750'
Set-TextValue $ws.Range("E47") '  -1.68%  '
Set-TextValue $ws.Range("D48") '2.033'
Set-TextValue $ws.Range("E48") '  -2.26%  '
Set-TextValue $ws.Range("D49") '123.89'
Set-TextValue $ws.Range("E49") '  -4.43%  '
Set-TextValue $ws.Range("D50") '1.215'
Set-TextValue $ws.Range("E50") '  -1.06%  '
Set-TextValue $ws.Range("D51") '0.07201'
Set-TextValue $ws.Range("E51") '  -1.78%  '
